$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "b_range" column (column I). This shifts "p_mutate" from J to I.
$ws.Columns.Item(9).Delete()

# 2. Bold the header row (A1:I1)
$ws.Range("A1:I1").Font.Bold = $true

# 3. Update max_iter for experiment 001 from 1000 to 100
$ws.Cells.Item(2, 2).Value = 100

# 4. Fill in experiment 002 (row 3)
$ws.Cells.Item(3, 1).Value = "'002"
$ws.Cells.Item(3, 2).Value = 500
$ws.Cells.Item(3, 3).Value = 0.8
$ws.Cells.Item(3, 4).Value = 100
$ws.Cells.Item(3, 5).Value = 100
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 6
$ws.Cells.Item(3, 8).Value = "'1,10"
$ws.Cells.Item(3, 9).Value = 0.3

# 5. Fill in experiment 003 (row 4)
$ws.Cells.Item(4, 1).Value = "'003"
$ws.Cells.Item(4, 2).Value = 1000
$ws.Cells.Item(4, 3).Value = 0.8
$ws.Cells.Item(4, 4).Value = 100
$ws.Cells.Item(4, 5).Value = 100
$ws.Cells.Item(4, 6).Value = 2
$ws.Cells.Item(4, 7).Value = 6
$ws.Cells.Item(4, 8).Value = "'1,10"
$ws.Cells.Item(4, 9).Value = 0.3

# 6. Add experiment 004 (row 5)
$ws.Cells.Item(5, 1).Value = "'004"
$ws.Cells.Item(5, 2).Value = 100
$ws.Cells.Item(5, 3).Value = 0.8
$ws.Cells.Item(5, 4).Value = 500
$ws.Cells.Item(5, 5).Value = 100
$ws.Cells.Item(5, 6).Value = 2
$ws.Cells.Item(5, 7).Value = 6
$ws.Cells.Item(5, 8).Value = "'1,10"
$ws.Cells.Item(5, 9).Value = 0.3

# 7. Add experiment 005 (row 6)
$ws.Cells.Item(6, 1).Value = "'005"
$ws.Cells.Item(6, 2).Value = 100
$ws.Cells.Item(6, 3).Value = 0.8
$ws.Cells.Item(6, 4).Value = 1000
$ws.Cells.Item(6, 5).Value = 100
$ws.Cells.Item(6, 6).Value = 2
$ws.Cells.Item(6, 7).Value = 6
$ws.Cells.Item(6, 8).Value = "'1,10"
$ws.Cells.Item(6, 9).Value = 0.3

# 8. Add experiment 006 (row 7)
$ws.Cells.Item(7, 1).Value = "'006"
$ws.Cells.Item(7, 2).Value = 100
$ws.Cells.Item(7, 3).Value = 0.8
$ws.Cells.Item(7, 4).Value = 1000
$ws.Cells.Item(7, 5).Value = 100
$ws.Cells.Item(7, 6).Value = 2
$ws.Cells.Item(7, 7).Value = 6
$ws.Cells.Item(7, 8).Value = "'1,10"
$ws.Cells.Item(7, 9).Value = 0.1

# 9. Add experiment 007 (row 8)
$ws.Cells.Item(8, 1).Value = "'007"
$ws.Cells.Item(8, 2).Value = 100
$ws.Cells.Item(8, 3).Value = 0.8
$ws.Cells.Item(8, 4).Value = 1000
$ws.Cells.Item(8, 5).Value = 100
$ws.Cells.Item(8, 6).Value = 2
$ws.Cells.Item(8, 7).Value = 6
$ws.Cells.Item(8, 8).Value = "'1,10"
$ws.Cells.Item(8, 9).Value = 0.5

# 10. Update the selected cell
[void]$ws.Range("F12").Select()

# 11. Match page orientation (portrait) as set in the source workbook
$ws.PageSetup.Orientation = 1
